# Daily attendance processing - 2026-01-27 09:57:52
# Normalize the "Recorded By" (column G) entries so the signed-in user's
# email is listed before the generic "System" actor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($current -eq $oldValue) {
        $cell.Value = $newValue
    }
}
